# Generate Report for Archive
#
# Refreshes the localization-status snapshot: the two files that have moved
# on from "Ready for handoff" into active translation
# (27e1c287-b513-49e7-b847-27c4d99ad29f and 7e7d8541-3990-4077-a7b6-09a55f7c5d0a)
# now report a status of "In Translation" everywhere that status is shown -
# the per-language detail sheets as well as the roll-up Overview sheet.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Per-language detail sheets: Status lives in column C, rows 3 and 4 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C3").Value = $newStatus
    $ws.Range("C4").Value = $newStatus
}

# --- Overview sheet: zh-cn status in column B, de-de status in column C ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus
